$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin border all around, centered horizontal, top vertical alignment
$cell = $ws.Range("B1")
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108  # xlCenter
$cell.VerticalAlignment = -4160    # xlTop
$cell.Borders.LineStyle = 1        # xlContinuous
$cell.Borders.Weight = 2           # xlThin

# Apply the same formatting to A2 by copying the formats from B1
$cell.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
